$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so values like "562.19" or "59.509.60"
# are stored as strings (matching the original inline-string cells) rather than
# being auto-converted to numbers by Excel.
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Range("D2").Value = '59.509.60'
$ws.Range("E2").Value = '  -2.28%  '
$ws.Range("D3").Value = '2.590.42'
$ws.Range("E3").Value = '  -2.10%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '562.19'
$ws.Range("E5").Value = '  -1.08%  '
$ws.Range("D6").Value = '143.23'
$ws.Range("E6").Value = '  -2.62%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.598'
$ws.Range("E8").Value = '  -1.78%  '
$ws.Range("D9").Value = '2.602.68'
$ws.Range("E9").Value = '  -2.65%  '
$ws.Range("E10").Value = '  -2.72%  '
$ws.Range("D11").Value = '0.105'
$ws.Range("E11").Value = '  -0.69%  '
$ws.Range("E12").Value = '  +10.06%  '
$ws.Range("E13").Value = '  +4.23%  '
$ws.Range("D14").Value = '3.043.93'
$ws.Range("E14").Value = '  -3.12%  '
$ws.Range("D15").Value = '23.29'
$ws.Range("E15").Value = '  +6.09%  '
$ws.Range("D16").Value = '59.432.41'
$ws.Range("E16").Value = '  -2.26%  '
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("D18").Value = '2.629.54'
$ws.Range("E18").Value = '  -1.46%  '
$ws.Range("D19").Value = '4.58'
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("D20").Value = '339.18'
$ws.Range("E20").Value = '  -1.44%  '
$ws.Range("D21").Value = '10.42'
$ws.Range("E21").Value = '  -0.57%  '
$ws.Range("D22").Value = '6.53'
$ws.Range("E22").Value = '  +2.45%  '
$ws.Range("E23").Value = '  +0.26%  '
$ws.Range("D24").Value = '63.42'
$ws.Range("E24").Value = '  -4.96%  '
$ws.Range("E25").Value = '  +7.30%  '
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("E27").Value = '  -1.75%  '
$ws.Range("E28").Value = '  +0.82%  '
$ws.Range("D29").Value = '0.0₃0781'
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("D31").Value = '6.21'
$ws.Range("E31").Value = '  -1.33%  '
$ws.Range("E32").Value = '  -2.30%  '
$ws.Range("D33").Value = '158.36'
$ws.Range("E33").Value = '  +1.64%  '
$ws.Range("D34").Value = '19.10'
$ws.Range("E34").Value = '  -0.87%  '
$ws.Range("E35").Value = '  -0.67%  '
$ws.Range("D36").Value = '1.17'
$ws.Range("E36").Value = '  +1.08%  '
$ws.Range("D37").Value = '0.902'
$ws.Range("E37").Value = '  -0.89%  '
$ws.Range("D38").Value = '0.872'
$ws.Range("E38").Value = '  -3.89%  '
$ws.Range("D39").Value = '37.43'
$ws.Range("E39").Value = '  -0.44%  '
$ws.Range("E40").Value = '  -1.95%  '
$ws.Range("D41").Value = '295.33'
$ws.Range("E41").Value = '  -2.90%  '
$ws.Range("E42").Value = '  +0.68%  '
$ws.Range("D43").Value = '139.77'
$ws.Range("E43").Value = '  +8.65%  '
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("E45").Value = '  -0.86%  '
$ws.Range("E46").Value = '  -1.57%  '
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").Value = '0.0533'
$ws.Range("E48").Value = '  -2.71%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = '0.0236'
$ws.Range("E49").Value = '  -0.29%  '
$ws.Range("D50").Value = '18.80'
$ws.Range("E50").Value = '  -0.81%  '
$ws.Range("D51").Value = '1.968.08'
$ws.Range("E51").Value = '  -0.01%  '

# Reset column D style back to the default/unstyled state (no explicit style index),
# matching the original workbook formatting.
$colD.Style = "Normal"

